$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$quote = [char]39

for ($row = 3; $row -le 26; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # column H - PERIOD TO EXPIRE
    $hCell.Value = $hCell.Value2 - 1

    $iCell = $ws.Cells.Item($row, 9)   # column I - LAST UPDATE
    # Leading apostrophe forces literal text so "04-Nov-2025" is not
    # auto-converted into a date serial value (matches the source file,
    # which stores this column as plain text).
    $iCell.Value = $quote + "04-Nov-2025"
}
